$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs / Agrp / Mc4r / MuSCs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Agrp"
$ws.Range("C2").Value = "Mc4r"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7829203333333333
$ws.Range("H2").Value = 2.348761
$ws.Range("I2").Value = 0.342817505175991
$ws.Range("J2").Value = 0.342817505175991
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.309635
$ws.Range("N2").Value = 0.928905
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.2424195374116667
$ws.Range("R2").Value = 2.181775836705
$ws.Range("S2").Value = 0.342817505175991
$ws.Range("T2").Value = 0.342817505175991

# Row 3: Inflammatory-Mac / Agrp / Mc4r / MuSCs
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Agrp"
$ws.Range("C3").Value = "Mc4r"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8792006666666667
$ws.Range("H3").Value = 2.637602
$ws.Range("I3").Value = 0.3849757967231252
$ws.Range("J3").Value = 0.3849757967231252
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.309635
$ws.Range("N3").Value = 0.928905
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.2722312984233333
$ws.Range("R3").Value = 2.45008168581
$ws.Range("S3").Value = 0.3849757967231252
$ws.Range("T3").Value = 0.3849757967231252

# Row 4: MuSCs / Agrp / Mc4r / MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Agrp"
$ws.Range("C4").Value = "Mc4r"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05636633333333333
$ws.Range("H4").Value = 0.169099
$ws.Range("I4").Value = 0.02468113925076025
$ws.Range("J4").Value = 0.02468113925076025
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.309635
$ws.Range("N4").Value = 0.928905
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.01745298962166666
$ws.Range("R4").Value = 0.157076906595
$ws.Range("S4").Value = 0.02468113925076025
$ws.Range("T4").Value = 0.02468113925076025

# Row 5: Resolving-Mac / Agrp / Mc4r / MuSCs
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Agrp"
$ws.Range("C5").Value = "Mc4r"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5652943333333332
$ws.Range("H5").Value = 1.695883
$ws.Range("I5").Value = 0.2475255588501236
$ws.Range("J5").Value = 0.2475255588501236
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.309635
$ws.Range("N5").Value = 0.928905
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.1750349109016666
$ws.Range("R5").Value = 1.575314198115
$ws.Range("S5").Value = 0.2475255588501236
$ws.Range("T5").Value = 0.2475255588501236
